# Insert a new weekly record at row 177 (pushing existing rows 177-238 down
# to 178-239) on the "Ajo" (garlic) price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 177; Excel shifts rows 177-238
# down to 178-239 and copies formatting (e.g. the date number format) from
# the row above into the new row.
$ws.Rows("177:177").Insert()

# Populate the newly inserted row 177 with the new weekly data point.
$ws.Range("A177").Value = 11
$ws.Range("B177").Value = "Vega Monumental Concepción"
$ws.Range("C177").Value = "Bíobío"
$ws.Range("D177").Value = 44924
$ws.Range("E177").Value = 8
$ws.Range("F177").Value = 100112003
$ws.Range("G177").Value = "Ajo"
$ws.Range("H177").Value = "Chino"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 400
$ws.Range("K177").Value = 13000
$ws.Range("L177").Value = 13500
$ws.Range("M177").Value = 13250
$ws.Range("N177").Value = '$/caja 10 kilos'
$ws.Range("O177").Value = "China"
$ws.Range("P177").Value = 1325
$ws.Range("Q177").Value = 10
$ws.Range("R177").Value = "Hortaliza"
